$wb = $excel.ActiveWorkbook

# 1) Status text update: "Ready for handoff" -> "In Translation"
#    (Overview!E2:F4, zh-cn!C2:C4, de-de!C2:C4 all hold this shared string.)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2:F4").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2:C4").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2:C4").Value = "In Translation"

# 2) Narrow the columns that used to be sized for "Ready for handoff" text
#    now that the shorter "In Translation" text is used.
$wsOverview.Columns.Item(5).ColumnWidth = 12.6
$wsOverview.Columns.Item(6).ColumnWidth = 12.6

$wsZhCn.Columns.Item(3).ColumnWidth = 12.6

$wsDeDe.Columns.Item(3).ColumnWidth = 12.6
